$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 becomes "atendido" with the values formerly in row 3
$ws.Range("A2").Value = "atendido"
$ws.Range("B2").Value = 667
$ws.Range("C2").Value = 82.95999999999999

# Row 3 becomes "violado" with the values formerly in row 2
$ws.Range("A3").Value = "violado"
$ws.Range("B3").Value = 137
$ws.Range("C3").Value = 17.04
